# "Invalid Product handler, status update, cart clearing"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Status column (B) -> flag rows 2-4 with the invalid-product handler message
$ws.Range("B2").Value = "Invalid Product exist"
$ws.Range("B3").Value = "Invalid Product exist"
$ws.Range("B4").Value = "Invalid Product exist"

# Status update: record the product that triggered the invalid-product check
$ws.Range("M2").Value = "Apple Magic Key Board"
$ws.Range("L3").Value = "Apple Magic Key Board"
$ws.Range("L4").Value = "Apple Magic Key Board"

# Cart clearing for row 5: remove the old cart items and replace with the new product
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("L5").Value = "EVOLVE 65 WIRELESS BLUETOOTH SINGLE EAR HEADSET"

# Update the active selection to reflect the last edited cell
$ws.Range("L5").Select()
